# Refresh the cryptos price table (Coin / Link / Price / Volume(1h)) with the
# latest scrape. Rows keep their position (A is a static 0-based index) but
# the coin that lands on a given row, along with its link/price/volume, can
# change from run to run as the ranking shifts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='67.863.52' }
    @{ Row=3; D='3.835.27'; E='  +1.74%  ' }
    @{ Row=4; E='  +0.00%  ' }
    @{ Row=5; D='603.57'; E='  +1.56%  ' }
    @{ Row=6; D='166.78'; E='  -0.14%  ' }
    @{ Row=7; D='1.00'; E='  -0.09%  ' }
    @{ Row=8; E='  +0.05%  ' }
    @{ Row=9; D='0.161'; E='  +1.06%  ' }
    @{ Row=10; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.452'; E='  +1.18%  ' }
    @{ Row=11; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='6.30'; E='  -0.52%  ' }
    @{ Row=12; E='  -0.24%  ' }
    @{ Row=13; D='35.88'; E='  -0.54%  ' }
    @{ Row=14; D='4.472.70'; E='  +1.77%  ' }
    @{ Row=15; D='3.844.39'; E='  +2.25%  ' }
    @{ Row=16; D='67.907.39'; E='  +0.46%  ' }
    @{ Row=17; D='18.44'; E='  +0.73%  ' }
    @{ Row=18; D='7.07'; E='  +1.27%  ' }
    @{ Row=19; E='  +0.55%  ' }
    @{ Row=20; D='465.82'; E='  +2.17%  ' }
    @{ Row=21; D='9.91'; E='  -0.63%  ' }
    @{ Row=22; D='0.701'; E='  +0.77%  ' }
    @{ Row=23; E='  -2.44%  ' }
    @{ Row=24; D='83.36'; E='  +0.26%  ' }
    @{ Row=25; D='12.10'; E='  +1.71%  ' }
    @{ Row=26; D='2.12'; E='  -0.36%  ' }
    @{ Row=27; D='10.08'; E='  -0.18%  ' }
    @{ Row=28; E='  -0.04%  ' }
    @{ Row=29; D='3.985.40'; E='  +1.90%  ' }
    @{ Row=30; D='2.79'; E='  +0.86%  ' }
    @{ Row=31; D='7.41'; E='  +2.02%  ' }
    @{ Row=32; D='2.23'; E='  +0.74%  ' }
    @{ Row=33; D='29.67'; E='  +0.32%  ' }
    @{ Row=34; B='RenzoRestakedETH'; C='https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'; D='3.779.07'; E='  +1.53%  ' }
    @{ Row=35; D='9.12'; E='  -0.17%  ' }
    @{ Row=36; B='Binance-PegBSC-USD'; C='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; D='1.00'; E='  +0.11%  ' }
    @{ Row=37; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.100'; E='  +0.07%  ' }
    @{ Row=38; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='3.30'; E='  +0.88%  ' }
    @{ Row=39; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.138'; E='  +0.34%  ' }
    @{ Row=40; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='1.00'; E='  +0.75%  ' }
    @{ Row=41; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='5.81'; E='  +1.41%  ' }
    @{ Row=42; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='0.999'; E='  +0.04%  ' }
    @{ Row=43; B='USDe'; C='https://coinranking.com/coin/exbfr2U-0+usde-usde'; D='1.00'; E='  +0.01%  ' }
    @{ Row=44; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='29.33'; E='  +8.95%  ' }
    @{ Row=45; B='Arweave'; C='https://coinranking.com/coin/7XWg41D1+arweave-ar'; D='44.79'; E='  -1.97%  ' }
    @{ Row=46; B='ONDO'; C='https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'; D='1.44'; E='  +15.51%  ' }
    @{ Row=47; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='47.84'; E='  -1.17%  ' }
    @{ Row=48; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.301'; E='  +0.78%  ' }
    @{ Row=49; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='149.60'; E='  +0.92%  ' }
    @{ Row=50; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.34'; E='  +0.41%  ' }
    @{ Row=51; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.84'; E='  +1.43%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B$r").Value = $u["B"] }
    if ($u.ContainsKey("C")) { $ws.Range("C$r").Value = $u["C"] }
    if ($u.ContainsKey("D")) {
        # Force text formatting so values like "1.00" / "0.100" keep their
        # exact printed form instead of being coerced into numbers.
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $u["D"]
    }
    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u["E"] }
}

Write-Host "Updated $($updates.Count) rows"